$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 122; this shifts all rows 122..203 down to 123..204
$ws.Rows.Item(122).Insert()

# Fill in the new row 122 with its data (columns A,B,C,E,F,G,H,I,R are constant
# across this data block, so copy them from the row below, which is the former
# row 122 now shifted to row 123)
$ws.Range("A122").Value = $ws.Range("A123").Value2
$ws.Range("B122").Value = $ws.Range("B123").Value2
$ws.Range("C122").Value = $ws.Range("C123").Value2
$ws.Range("D122").Value = 44603
$ws.Range("E122").Value = $ws.Range("E123").Value2
$ws.Range("F122").Value = $ws.Range("F123").Value2
$ws.Range("G122").Value = $ws.Range("G123").Value2
$ws.Range("H122").Value = $ws.Range("H123").Value2
$ws.Range("I122").Value = $ws.Range("I123").Value2
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 9500
$ws.Range("L122").Value = 10000
$ws.Range("M122").Value = 9750
$ws.Range("N122").Value = '$/caja 80 unidades'
$ws.Range("O122").Value = 'Región del Maule'
$ws.Range("P122").Value = 122
$ws.Range("Q122").Value = 80
$ws.Range("R122").Value = $ws.Range("R123").Value2
